$d = $word.ActiveDocument
$d.Content.Find.Execute("GPA: 3.5", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GPA: 3.6", 2)
